$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: the timestamp "SAT Dec 23 12:39:06 PST 2017" used to be split across
# two runs ("SAT Dec 23" + " 12:39:06 PST 2017"). Re-run a Find/Replace over
# the same text so Word collapses it back into a single run/text node.
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "SAT Dec 23 12:39:06 PST 2017", $false, $false, $false, $false, $false,
    $true, 1, $false, "SAT Dec 23 12:39:06 PST 2017", 2)
if (-not $found1) {
    throw "Could not find the 'SAT Dec 23 12:39:06 PST 2017' timestamp to normalize."
}

# ---------------------------------------------------------------------------
# Edit 2: append a brand new "MON Dec 25" purchase record right after the
# existing "Amount balance ... - 160869.0" paragraph (the end of the SAT Dec
# 23 record), mirroring the layout used by every other purchase block in the
# document.
# ---------------------------------------------------------------------------
$r = $d.Content
$found2 = $r.Find.Execute("160869.0", $false, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the '- 160869.0' anchor to insert the new record after."
}

# Collapse to the end of the match, then expand/collapse to land exactly on
# the paragraph mark that ends the "Amount balance" paragraph.
$r.Collapse(0)
$r.Expand(4)
$r.Collapse(0)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$font = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>'

function PlainPara([string]$innerRuns, [bool]$bold = $false) {
    if ($bold) {
        $rpr = "<w:rPr>$font<w:b/></w:rPr>"
    } else {
        $rpr = "<w:rPr>$font</w:rPr>"
    }
    return "<w:p><w:pPr><w:pStyle w:val=`"PlainText`"/>$rpr</w:pPr>$innerRuns</w:p>"
}

function LabeledRun([string]$label, [int]$tabCount, [string]$value, [bool]$bold = $false) {
    if ($bold) {
        $rpr = "<w:rPr>$font<w:b/></w:rPr>"
    } else {
        $rpr = "<w:rPr>$font</w:rPr>"
    }
    $runs = "<w:r>$rpr<w:t>$label</w:t></w:r>"
    for ($i = 0; $i -lt $tabCount; $i++) {
        $runs += "<w:r>$rpr<w:tab/></w:r>"
    }
    $runs += "<w:r>$rpr<w:tab/><w:t>$value</w:t></w:r>"
    return $runs
}

$emptyBold = PlainPara "" $true
$timestamp = PlainPara ("<w:r><w:rPr>$font</w:rPr><w:t>MON Dec 25</w:t></w:r>" + `
    "<w:r><w:rPr>$font</w:rPr><w:t xml:space=`"preserve`"> 12:10:45 PST 2017</w:t></w:r>")
$personName = PlainPara (LabeledRun "Person Name" 3 "- HHN")
$billNumber = PlainPara (LabeledRun "Bill number" 3 "- 746")
$dashes = PlainPara "<w:r><w:rPr>$font</w:rPr><w:t>---------------------------------------------------------------</w:t></w:r>"
$itemName = PlainPara (LabeledRun "Item Name" 3 "- CARROT2")
$pockets = PlainPara (LabeledRun "Number of Pockets" 2 "- 1")
$kgs = PlainPara (LabeledRun "Number of KGs" 2 "- 93")
$rate = PlainPara (LabeledRun "Rate" 4 "- 48")
$transport = PlainPara ("<w:r><w:rPr>$font</w:rPr><w:t>Transport &amp; Miscellaneous</w:t></w:r>" + `
    "<w:r><w:rPr>$font</w:rPr><w:tab/><w:t>- 15</w:t></w:r>")
$totalPrice = PlainPara (LabeledRun "Total Price" 3 "- 4479.0")
$amountBalance = PlainPara (LabeledRun "Amount balance" 2 "- 165348.0" $true) $true
$trailingEmpty = PlainPara ""
$trailingEmptyBold = PlainPara "" $true

$xml = "<w:p $w>" + $emptyBold.Substring(4) + $timestamp + $personName + $billNumber + `
    $dashes + $itemName + $pockets + $kgs + $rate + $transport + $totalPrice + `
    $amountBalance + $trailingEmpty + $trailingEmptyBold

$r.InsertXML($xml)
